$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D14").Value = 27600
$ws.Range("D43").Value = 218600
$ws.Range("D45").Value = 63600
$ws.Range("D46").Value = 1155200
$ws.Range("D48").Value = 341000
$ws.Range("D49").Value = 907700
$ws.Range("D52").Value = 129700
$ws.Range("D54").Value = 2026800
$ws.Range("D58").Value = 551700
$ws.Range("D59").Value = 228900
$ws.Range("D60").Value = 553200
$ws.Range("D62").Value = 154000
$ws.Range("D66").Value = 1291900
$ws.Range("D72").Value = -69378600
$ws.Range("D76").Value = 734900
